# update version name & version code.
#
# Bumps the APK version strings referenced by the "安装文件" (install file)
# column (I/J) for both adaptation rows:
#   - 企鹅极光1S (PENGUIN1S A4062) row: 1.15 -> 1.18
#   - TCL Android TV (55A261)      row: 1.15 -> 1.19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 企鹅极光1S
$ws.Range("I2").Value = "CQLive-PENGUIN1SA4062-1.18"
$ws.Range("J2").Value = "CQLive-PENGUIN1SA4062-1.18.apk"

# Row 3: TCL Android TV / 55A261
$ws.Range("I3").Value = "CQLive-TCL55A261-1.19"
$ws.Range("J3").Value = "CQLive-TCL55A261-1.19.apk"

# Leave the selection where the author last left it when saving.
$ws.Range("J9").Select()
